$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$a2 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f45a9b8a370>),
                ('model',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=1.0, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.2, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=7, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])
'@
$ws.Cells.Item(2,1).Value = $a2

$ws.Cells.Item(2,2).Value = 0.6666666666666666
$c2 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f45a9ace8b0>, 'scaler': None, 'model__subsample': 0.5, 'model__n_estimators': 50, 'model__max_depth': 7, 'model__learning_rate': 0.01, 'model__gamma': 0.2, 'model__colsample_bytree': 1.0}
'@
$ws.Cells.Item(2,3).Value = $c2
$ws.Cells.Item(2,4).Value = 0.3333333333333333
$ws.Cells.Item(2,5).Value = '[1 1 0 0 1 0 0 0 0 1 0 1]'
$ws.Cells.Item(2,6).Value = '[0 1 1 0 0 1 1 1 1 1 0 0]'
$ws.Cells.Item(2,7).Value = 77
$ws.Cells.Item(2,8).Value = 0.9669279907084787
$ws.Cells.Item(2,9).Value = 0.007192539814404358
$ws.Cells.Item(2,10).Value = 0.5635307781649244
$ws.Cells.Item(2,11).Value = 0.06714518070469799

# Row 3
$a3 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f45043fa3a0>),
                ('model',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.1, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=7, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=100,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])
'@
$ws.Cells.Item(3,1).Value = $a3

$ws.Cells.Item(3,2).Value = 0.638095238095238
$c3 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3adc0c47f0>, 'scaler': None, 'model__subsample': 0.5, 'model__n_estimators': 100, 'model__max_depth': 7, 'model__learning_rate': 0.01, 'model__gamma': 0.1, 'model__colsample_bytree': 0.5}
'@
$ws.Cells.Item(3,3).Value = $c3
$ws.Cells.Item(3,4).Value = 0.7368421052631579
$ws.Cells.Item(3,5).Value = '[1 1 0 1 0 0 1 0 1 1 1 0]'
$ws.Cells.Item(3,6).Value = '[1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Cells.Item(3,7).Value = 69
$ws.Cells.Item(3,8).Value = 0.9816190476190475
$ws.Cells.Item(3,9).Value = 0.005665436470911159
$ws.Cells.Item(3,10).Value = 0.5132698412698411
$ws.Cells.Item(3,11).Value = 0.08949786041996534

# Row 4
$a4 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3adc0c4eb0>),
                ('model',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.1, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.1,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=7, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])
'@
$ws.Cells.Item(4,1).Value = $a4

$ws.Cells.Item(4,2).Value = 0.6095238095238095
$c4 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f4504368220>, 'scaler': None, 'model__subsample': 0.8, 'model__n_estimators': 50, 'model__max_depth': 7, 'model__learning_rate': 0.1, 'model__gamma': 0.1, 'model__colsample_bytree': 0.5}
'@
$ws.Cells.Item(4,3).Value = $c4
$ws.Cells.Item(4,4).Value = 0.7058823529411765
$ws.Cells.Item(4,5).Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Cells.Item(4,6).Value = '[1 1 1 1 1 1 0 0 1 1 1 0]'
$ws.Cells.Item(4,7).Value = 42
$ws.Cells.Item(4,8).Value = 0.979182754182754
$ws.Cells.Item(4,9).Value = 0.007740351256799776
$ws.Cells.Item(4,10).Value = 0.5241956241956242
$ws.Cells.Item(4,11).Value = 0.09831409418318524

# Row 5
$a5 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a6c5373d0>),
                ('model',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.8, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.1,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=3, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=100,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])
'@
$ws.Cells.Item(5,1).Value = $a5

$ws.Cells.Item(5,2).Value = 0.6190476190476191
$c5 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c439940>, 'scaler': None, 'model__subsample': 0.8, 'model__n_estimators': 100, 'model__max_depth': 3, 'model__learning_rate': 0.1, 'model__gamma': 0, 'model__colsample_bytree': 0.8}
'@
$ws.Cells.Item(5,3).Value = $c5
$ws.Cells.Item(5,4).Value = 0.6666666666666666
$ws.Cells.Item(5,5).Value = '[1 1 0 0 0 0 1 0 1 1 1 1]'
$ws.Cells.Item(5,6).Value = '[1 0 0 1 0 1 1 1 1 0 1 1]'
$ws.Cells.Item(5,7).Value = 11
$ws.Cells.Item(5,8).Value = 0.9850732600732603
$ws.Cells.Item(5,9).Value = 0.004934424748355733
$ws.Cells.Item(5,10).Value = 0.5063492063492063
$ws.Cells.Item(5,11).Value = 0.07077801350542343

# Row 6
$a6 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a6c439100>),
                ('model',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0.1, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.1,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=3, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])
'@
$ws.Cells.Item(6,1).Value = $a6

$ws.Cells.Item(6,2).Value = 0.5904761904761905
$c6 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c3f96d0>, 'scaler': None, 'model__subsample': 0.5, 'model__n_estimators': 50, 'model__max_depth': 3, 'model__learning_rate': 0.1, 'model__gamma': 0.1, 'model__colsample_bytree': 0.5}
'@
$ws.Cells.Item(6,3).Value = $c6
$ws.Cells.Item(6,4).Value = 0.7692307692307692
$ws.Cells.Item(6,5).Value = '[1 1 1 1 0 0 0 0 1 1 0 0]'
$ws.Cells.Item(6,6).Value = '[1 1 1 0 0 1 0 1 1 1 0 0]'
$ws.Cells.Item(6,7).Value = 14
$ws.Cells.Item(6,8).Value = 0.9839517625231912
$ws.Cells.Item(6,9).Value = 0.004337687213876009
$ws.Cells.Item(6,10).Value = 0.5083487940630798
$ws.Cells.Item(6,11).Value = 0.09994895771068718

